$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp shown in the title cell (A1).
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 09:56"

# --- Country reordering (new country inserted just above an existing one,
#     pushing the existing country's unchanged data down one row) ---

# Armenia now sorts before Argelia (rows 63/64): Armenia gets fresh stats,
# Argelia keeps its previous stats but moves down to row 64.
$ws.Range("A63").Value = "Armenia"
$ws.Range("B63").Value = 52496
$ws.Range("C63").Value = 571
$ws.Range("D63").Value = 44672
$ws.Range("E63").Value = 6847
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 977

$ws.Range("A64").Value = "Argelia"
$ws.Range("B64").Value = 51995
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 36482
$ws.Range("E64").Value = 13757
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1756

# Georgia now sorts before Zimbabue (rows 113/114): Georgia gets fresh
# stats, Zimbabue keeps its previous stats but moves down to row 114.
$ws.Range("A113").Value = "Georgia"
$ws.Range("B113").Value = 8118
$ws.Range("C113").Value = 554
$ws.Range("D113").Value = 4244
$ws.Range("E113").Value = 3824
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 50

$ws.Range("A114").Value = "Zimbabue"
$ws.Range("B114").Value = 7885
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 6327
$ws.Range("E114").Value = 1330
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 228

# Letonia now sorts before Yemen (rows 157/158): Letonia gets fresh
# stats, Yemen keeps its previous stats but moves down to row 158.
$ws.Range("A157").Value = "Letonia"
$ws.Range("B157").Value = 2086
$ws.Range("C157").Value = 67
$ws.Range("D157").Value = 1307
$ws.Range("E157").Value = 741
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 38

$ws.Range("A158").Value = "Yemen"
$ws.Range("B158").Value = 2041
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 1320
$ws.Range("E158").Value = 132
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 589

# --- Plain statistic refreshes (country / row unchanged) ---

# Row 27: Israel
$ws.Range("B27").Value = 264857
$ws.Range("C27").Value = 414
$ws.Range("D27").Value = 193002
$ws.Range("E27").Value = 70173

# Row 60
$ws.Range("B60").Value = 57812
$ws.Range("C60").Value = 12
$ws.Range("E60").Value = 223

# Row 77
$ws.Range("B77").Value = 30575
$ws.Range("C77").Value = 858
$ws.Range("D77").Value = 7470
$ws.Range("E77").Value = 22283
$ws.Range("G77").Value = 10
$ws.Range("H77").Value = 822

# Row 81
$ws.Range("B81").Value = 27136
$ws.Range("C81").Value = 15
$ws.Range("D81").Value = 24866
$ws.Range("E81").Value = 1376
